$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, since the source data keeps these as literal strings
# (inline strings) rather than numeric cells.

$ws.Range("D2").Value = "60.967.97"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.385.74"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.41"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.93"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.64"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "3.965.15"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.71"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "3.392.17"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "61.087.73"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  -3.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.64"
$ws.Range("E19").Value = "  -5.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.94"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.22"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.87"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  -5.12%  "
$ws.Range("D26").Value = "3.521.24"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.31"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.98"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.20"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.02"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.416.15"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0768"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.97"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "2.451.43"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.02"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  +6.05%  "
